$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 256, pushing the existing rows 256-290 down to 257-291
# (this also grows the sheet dimension from A1:R290 to A1:R291).
$ws.Rows("256:256").Insert()

# Populate the newly inserted row 256 with this week's price entry.
$ws.Range("A256").Value = 6
$ws.Range("B256").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C256").Value = "Metropolitana"
$ws.Range("D256").Value = 44984
$ws.Range("E256").Value = 13
$ws.Range("F256").Value = 100112001
$ws.Range("G256").Value = "Berenjena"
$ws.Range("H256").Value = "Sin especificar"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 130
$ws.Range("K256").Value = 6000
$ws.Range("L256").Value = 6000
$ws.Range("M256").Value = 6000
$ws.Range("N256").Value = "$/caja 50 unidades"
$ws.Range("O256").Value = "Región de Arica y Parinacota"
$ws.Range("P256").Value = 120
$ws.Range("Q256").Value = 50
$ws.Range("R256").Value = "Hortaliza"
